# The title and author lines get "typed" back in word-by-word, each
# word (and each inter-word space) landing in its own run instead of
# being collapsed into a single run of text.
#
# Word's Range.Text / InsertBefore / TypeText APIs all splice text into
# the existing run when the surrounding formatting is identical, so the
# only reliable way to land distinct <w:r> siblings with matching (empty)
# formatting is to replace the whole paragraph (InsertXML REPLACES the
# target range) with freshly authored OOXML that already has the runs
# split out.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-ParagraphXml([string]$style, [string[]]$segments) {
    $runs = ($segments | ForEach-Object {
        '<w:r><w:t xml:space="preserve">' + $_ + '</w:t></w:r>'
    }) -join ''
    return '<w:p ' + $wNs + '><w:pPr><w:pStyle w:val="' + $style + '"/></w:pPr>' + $runs + '</w:p>'
}

$d = $word.ActiveDocument

$titleSegments = @("Sigma", " ", "Notation:", " ", "Answers")
$authorSegments = @("Ifan", " ", "Howells-Baines,", " ", "Mark", " ", "Toner")

$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertXML((New-ParagraphXml "Title" $titleSegments))

$authorPara = $d.Paragraphs(2)
$authorPara.Range.InsertXML((New-ParagraphXml "Author" $authorSegments))
